$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.749.69'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '2.816.10'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = "'352.53"
$ws.Range('E5').Value = '  +5.83%  '
$ws.Range('D6').Value = "'112.91"
$ws.Range('E6').Value = '  -3.14%  '
$ws.Range('E7').Value = '  +5.10%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +3.85%  '
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = "'19.88"
$ws.Range('D14').Value = "'7.74"
$ws.Range('E14').Value = '  +0.99%  '
$ws.Range('D15').Value = '3.261.11'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').Value = '2.823.81'
$ws.Range('E16').Value = '  +1.76%  '
$ws.Range('D17').Value = "'0.885"
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').Value = '51.707.56'
$ws.Range('E18').Value = '  +0.13%  '
$ws.Range('E19').Value = '  +9.10%  '
$ws.Range('D20').Value = "'3.18"
$ws.Range('E20').Value = '  -4.24%  '
$ws.Range('D21').Value = "'13.37"
$ws.Range('E21').Value = '  -1.32%  '
$ws.Range('D22').Value = '0.0₃0993'
$ws.Range('E22').Value = '  +1.76%  '
$ws.Range('D23').Value = "'270.35"
$ws.Range('E23').Value = '  -2.99%  '
$ws.Range('D24').Value = "'69.67"
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  +2.26%  '
$ws.Range('D26').Value = "'26.70"
$ws.Range('D27').Value = "'1.00"
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('E30').Value = '  -1.78%  '
$ws.Range('D31').Value = "'33.98"
$ws.Range('E31').Value = '  -3.08%  '
$ws.Range('D32').Value = "'50.57"
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('E33').Value = '  +26.95%  '
$ws.Range('E34').Value = '  +4.35%  '
$ws.Range('D35').Value = "'5.28"
$ws.Range('E35').Value = '  +4.88%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E37').Value = '  +0.05%  '
$ws.Range('E38').Value = '  -1.76%  '
$ws.Range('D39').Value = "'3.22"
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('E40').Value = '  -5.70%  '
$ws.Range('D41').Value = "'23.88"
$ws.Range('E41').Value = '  +3.11%  '
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('D43').Value = "'126.31"
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('D44').Value = "'2.51"
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('D46').Value = '2.078.84'
$ws.Range('E46').Value = '  -0.55%  '
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('E49').Value = '  +2.09%  '
$ws.Range('D50').Value = "'0.932"
$ws.Range('E50').Value = '  +6.41%  '
$ws.Range('D51').Value = "'60.72"
$ws.Range('E51').Value = '  +0.18%  '
